# Update torque reading values (Test 1-5) for each Allowance Range row so the
# live torque display reflects freshly emitted readings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "220.8 - 239.2" allowance range
$ws.Range("B2").Value = 223.9
$ws.Range("C2").Value = 239.2
$ws.Range("D2").Value = 235.4
$ws.Range("E2").Value = 237.3
$ws.Range("F2").Value = 237.2

# Row 3: "144.0 - 156.0" allowance range
$ws.Range("B3").Value = 153.8
$ws.Range("C3").Value = 153.8
$ws.Range("D3").Value = 152.8
$ws.Range("E3").Value = 154.4
$ws.Range("F3").Value = 154.3

# Row 4: "67.2 - 72.8" allowance range
$ws.Range("B4").Value = 68.90000000000001
$ws.Range("C4").Value = 68.2
$ws.Range("D4").Value = 69.09999999999999
$ws.Range("E4").Value = 69
$ws.Range("F4").Value = 69.3
